# Tidied up references to supplementary material for reconciliations
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mentioned_in_text")

# Update the "Name" column (A) for the three reconciliation rows to
# point to the supplementary-files-on-GitHub wording instead of the
# placeholder "Supplementary Figure X" text.
$ws.Range("A11").Value = "supplementary files with the full reconciliation for GRK on GitHub"
$ws.Range("A12").Value = "supplementary files with the full reconciliation for PLC on GitHub"
$ws.Range("A16").Value = "supplementary files with the full reconciliation for opsins on GitHub"

# Widen column A to fit the new, longer text (and drop the old "best fit"
# autosizing in favour of an explicit width).
$ws.Columns.Item(1).ColumnWidth = 31.81640625

# Update the selection shown when the sheet is reopened.
$ws.Range("A17").Select()
